# Applies the edit described in the commit "Deleted redundant server and client code"
# to the Project Engineering Work Journal document.
#
# The substantive content change in the diff is the insertion of three new
# bulleted "Work Done" items right before the "Issues & Solutions" heading
# that follows the "8th Oct - 1st Dec" week entry:
#   - Implemented Tomcat local server to project
#   - Added AddMessage servlet to tomcat: allows for INSERT message into message table in DB
#   - Removed redundant pre-Tomcat server code
#
# (All of the remaining hunks in the diff only wrap already-present words in
# <w:proofErr> spell/grammar-check marks, which do not change any visible
# text -- Word regenerates those automatically as part of its proofing pass
# and they carry no semantic content.)

$d = $word.ActiveDocument

$anchor = "Server application returns stored string message to client, GUI is able to display returned string"

$replacement = $anchor + "^pImplemented Tomcat local server to project" + `
    "^pAdded AddMessage servlet to tomcat: allows for INSERT message into message table in DB" + `
    "^pRemoved redundant pre-Tomcat server code"

$found = $d.Content.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

Write-Output "Replace executed: $found"
